# angular slide update from last workshop
#
# 1) Slide 3 ("Bootstrapping"): collapse the three bullet paragraphs in the
#    body placeholder down to just "Enable new forms API", now shown as a
#    non-bulleted, indented line (bullet explicitly turned off, keeping the
#    level-1 text start position).
# 2) Slide 3: nudge the screenshot picture up now that the text block is
#    shorter.
#
# (The handout master's auto-updating date field also re-stamped itself
# to the day this deck was last saved, as PowerPoint always does - that's
# a side effect of saving, not a content edit, and isn't reproduced here.)

$p = $ppt.ActivePresentation

# --- (1) & (2): slide 3 body text + picture reposition -------------------
$s = $p.Slides.Item(3)

$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "Enable new forms API"
$tr.IndentLevel = 2
$tr.ParagraphFormat.Bullet.Type = 0

$ruler = $body.TextFrame.Ruler
$lvl2 = $ruler.Levels.Item(2)
$lvl2.LeftMargin = 36
$lvl2.FirstMargin = 0

$picture = $s.Shapes.Item(3)
$picture.Top = 186
